# [silverfox] new npc temp-data entry
# Adds two new "cleaner" style rows (cid 5119 / 5120) to Sheet1 and
# highlights the newly entered cells with a purple fill + white text,
# matching the author's manual "new data" callout formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 18 : cleaner_style / cid 5119 --------------------------------
$ws.Range("A18").Value = 5119
$ws.Range("B18").Value = "cleaner_style"
$ws.Range("C18").Value = "클리너 스타일"
$ws.Range("D18").Value = "style"
$ws.Range("E18").Value = "Weakness"
$ws.Range("F18").Value = "WalkersWorkshop"
$ws.Range("G18").Value = "no107_skill"
$ws.Range("H18").Value = "{(0.5)}"
$ws.Range("I18").Value = "melee"
$ws.Range("K18").Value = 10

# ---- Row 19 : cleaner_R_style / cid 5120 -------------------------------
$ws.Range("A19").Value = 5120
$ws.Range("B19").Value = "cleaner_R_style"
$ws.Range("C19").Value = "깔끔한 클리너 스타일"
$ws.Range("D19").Value = "style"
$ws.Range("E19").Value = "Weakness"
$ws.Range("F19").Value = "SmithsSmithy"
$ws.Range("G19").Value = "doncina_skill01"
$ws.Range("H19").Value = "{(1.0)}"
$ws.Range("I19").Value = "range"
$ws.Range("K19").Value = 15

# ---- Highlight the newly-entered cells (purple fill / white text) -----
$row18Fill = $ws.Range("E18:K18")
$row18Fill.Interior.Color = 10498160
$row18Fill.Font.ThemeColor = 2

$row19Fill = $ws.Range("E19:K19")
$row19Fill.Interior.Color = 10498160
$row19Fill.Font.ThemeColor = 2

# ---- Move the selection below the newly-added rows ---------------------
[void]$ws.Range("A20").Select()
